# Add a new row of data (row 20) to Sheet1 and apply percentage formatting
# to the "Proportion" columns (C and E) across the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Append new data row 20 ---
$ws.Range("A20").Value = 7150
$ws.Range("B20").Value = 5.02
$ws.Range("C20").Value = 1.201
$ws.Range("D20").Value = 1.7224
$ws.Range("E20").Value = 1.1035
$ws.Range("G20").Value = 430

# New row's Mean_Duration columns (B, D) keep the same "0.0000" number
# format already used by the rows above them.
$ws.Range("B20").NumberFormat = "0.0000"
$ws.Range("D20").NumberFormat = "0.0000"

# --- Apply percentage number format (0.00%) to the Proportion columns ---
$ws.Range("C3:C20").NumberFormat = "0.00%"
$ws.Range("E3:E20").NumberFormat = "0.00%"

# --- Update selection to match the post-edit state ---
$ws.Range("M18").Select()
